# Append two new daily COVID-19 data rows (2020-06-29 and 2020-06-30) to the
# "Tabela1" table on the active sheet, then grow the table/selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 111 (2020-06-29 / serial 44011) -----------------------------------
# Clone the formatting of the current last data row (110) so the new row
# inherits the same number formats / fonts / borders, then overwrite the
# values that actually change. Column I keeps its "111*" footnote label and
# column J keeps its "0", matching the existing pattern, so they are left
# exactly as copied.
$ws.Range("A110:J110").Copy($ws.Range("A111:J111"))
$ws.Range("A111").Value = 44011
$ws.Range("B111").Value = 100330
$ws.Range("C111").Value = 1085
$ws.Range("D111").Value = 1600
$ws.Range("E111").Value = 15
$ws.Range("F111").Value = 8
$ws.Range("G111").Value = 0
$ws.Range("H111").Value = 0

# --- Row 112 (2020-06-30 / serial 44012) -----------------------------------
# Row 112 uses the "un-bordered" formatting variant (as in the source file),
# so clone it from an existing row that already carries that look.
$ws.Range("A100:J100").Copy($ws.Range("A112:J112"))
$ws.Range("A112").Value = 44012
$ws.Range("B112").Value = 101729
$ws.Range("C112").Value = 1399
$ws.Range("D112").Value = 1613
$ws.Range("E112").Value = 13
$ws.Range("F112").Value = 8
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 0
$ws.Range("I112").Value = "111*"
$ws.Range("J112").Value = 0

# --- Grow the table / autofilter to cover the two new rows -----------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:J112"))

# --- Match the workbook's last selection ------------------------------------
$ws.Range("J112").Select() | Out-Null

Write-Host "Added rows 111:112 and resized Tabela1 to A1:J112"
